# Daily attendance processing - 2025-11-22 20:45:29
# Swap the order of the comma-separated "Recorded By" entries in column G
# whenever the cell currently lists "dnasr281@gmail.com" first, e.g.
#   "dnasr281@gmail.com, System"        -> "System, dnasr281@gmail.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$target = "dnasr281@gmail.com"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -eq 2 -and $parts[0].Trim() -eq $target) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            $cell.Value = "$second, $first"
        }
    }
}
